$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.845.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +10.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.443.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +10.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "476.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +23.30%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.436.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.99%  "

# Row 10
$ws.Range("E10").Value = "  +14.91%  "

# Row 11
$ws.Range("E11").Value = "  +10.73%  "

# Row 12
$ws.Range("E12").Value = "  +12.62%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.123"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.863.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.987.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +14.64%  "

# Row 17
$ws.Range("E17").Value = "  +20.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.440.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.64%  "

# Row 19
$ws.Range("E19").Value = "  +13.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +20.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.42%  "

# Row 25
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.399"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +22.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.535.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.74%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.69%  "

# Row 30
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0759"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.02%  "

# Row 34
$ws.Range("E34").Value = "  +16.06%  "

# Row 35
$ws.Range("E35").Value = "  +14.26%  "

# Row 36
$ws.Range("E36").Value = "  +19.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.28%  "

# Row 38
$ws.Range("E38").Value = "  +14.80%  "

# Row 39
$ws.Range("E39").Value = "  +7.53%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41
$ws.Range("E41").Value = "  +14.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.596"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0541"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.40%  "

# Row 44
$ws.Range("E44").Value = "  +19.32%  "

# Row 45
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +28.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "254.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +38.31%  "

# Row 48
$ws.Range("E48").Value = "  +16.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.895.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.87%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.94%  "
